$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 gains a new weigh-in log entry in columns F/G (date + weight),
# matching the pattern used by the other rows in that mini-table
# (F column holds a date styled like F3/F6/F8/..., G holds the plain number).
# Copy the date cell format from an existing F-column entry (F3) so the new
# cell reuses the existing date style instead of registering a new one.
$ws.Range("F3").Copy()
$ws.Range("F13").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F13").Value = 43606        # 2019-05-21
$ws.Range("G13").Value = 89.8

# Update the active selection to reflect where the user ended up (F16).
$ws.Range("F16").Select()
